$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force D and E columns to be treated as text so numeric-looking values
# (e.g. "1.004", "326.82") are preserved as strings, matching source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.484.25'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '1.974.22'
$ws.Range('E3').Value = '  +3.90%  '
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '326.82'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '0.4664'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').Value = '0.3920'
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('D9').Value = '46.20'
$ws.Range('E9').Value = '  -1.03%  '
$ws.Range('D10').Value = '0.07957'
$ws.Range('E10').Value = '  +0.93%  '
$ws.Range('D11').Value = '0.9912'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '22.79'
$ws.Range('E12').Value = '  +3.83%  '
$ws.Range('D13').Value = '1.968.69'
$ws.Range('E13').Value = '  +3.47%  '
$ws.Range('D14').Value = '7.184'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Value = '5.839'
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').Value = '0.07090'
$ws.Range('E16').Value = '  +1.48%  '
$ws.Range('D17').Value = '87.67'
$ws.Range('E17').Value = '  -0.94%  '
$ws.Range('D18').Value = '1.005'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('E20').Value = '  +1.27%  '
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').Value = '29.479.51'
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('D23').Value = '0.5027'
$ws.Range('E23').Value = '  +4.42%  '
$ws.Range('D24').Value = '5.543'
$ws.Range('E24').Value = '  +4.23%  '
$ws.Range('D25').Value = '11.15'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('D26').Value = '2.202.37'
$ws.Range('E26').Value = '  +3.25%  '
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('D28').Value = '158.49'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').Value = '19.50'
$ws.Range('E29').Value = '  +0.41%  '
$ws.Range('D30').Value = '5.781'
$ws.Range('E30').Value = '  -3.86%  '
$ws.Range('D31').Value = '119.63'
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').Value = '0.09420'
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('D34').Value = '0.8937'
$ws.Range('E34').Value = '  -1.76%  '
$ws.Range('D35').Value = '5.234'
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('D36').Value = '1.324'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').Value = '3.169'
$ws.Range('E37').Value = '  -1.56%  '
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('E39').Value = '  -1.60%  '
$ws.Range('D40').Value = '0.02104'
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('D41').Value = '7.769'
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.5722'
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = '0.000003176'
$ws.Range('E43').Value = '  +51.77%  '
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').Value = '9.671'
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('D46').Value = '2.748'
$ws.Range('E46').Value = '  +6.49%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.5351'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '11.72'
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('D49').Value = '2.200'
$ws.Range('E49').Value = '  +1.19%  '
$ws.Range('D50').Value = '0.06925'
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.836'
$ws.Range('E51').Value = '  -0.91%  '
